$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.674.52"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.071.26"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.18"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.97"
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0781"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D12").Value = "2.377.04"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.75"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.80"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.769"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "2.056.36"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "37.615.28"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.07"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.48"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.29"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.38"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  -5.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.82"
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.37"
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.32"
$ws.Range("E39").Value = "  -5.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0973"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.75"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "1.453.27"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.16"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.42"
$ws.Range("E46").Value = "  +5.08%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.05"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.01"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "2.261.11"
$ws.Range("E51").Value = "  -1.75%  "
